# Update the "# of tomatoes" counts (column K) on the "Sample inventory" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample inventory")

$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 4
$ws.Range("K8").Value = 4
$ws.Range("K10").Value = 4
$ws.Range("K11").Value = 9

# Reflect the last active cell selection recorded in the saved file
$ws.Range("J13").Select()

$wb.Save()
